$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "No" column (column A). This shifts NIM -> A and Nama -> B,
# and the old column C becomes empty/removed.
$ws.Range("A:A").Delete()

# Header row correction: A1 should read "NIM" (already true after shift),
# and B1 should read "Nama" (already true after shift). Nothing else needed there.

# Row 4 (NIM 71231014) needs to be stored as text, not a number (auto correction).
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "71231014"
$ws.Cells.Item(4, 1).Style = "Normal"
